$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.611.79"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -0.91%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.532.41"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -1.08%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.01%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'309.54"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -1.38%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'100.22"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +3.72%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.570"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -1.08%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.11%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.529"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -2.02%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'35.99"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +1.54%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  -0.82%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'7.34"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -1.14%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  +0.37%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'2.925.25"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.97%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'15.95"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +6.07%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'2.579.55"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +2.30%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'0.821"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -2.32%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'42.599.13"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.95%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'6.83"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.01%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.46%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'12.22"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -2.27%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'69.24"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.23%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'243.33"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -3.42%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'2.90"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -1.42%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'2.05"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -0.87%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  +0.04%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'25.94"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -2.80%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'2.33"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -3.97%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'39.37"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -1.44%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'10.14"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -0.34%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'156.36"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.21%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  -0.65%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  +13.77%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.0797"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -1.04%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'2.63"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -2.53%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'2.03"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -3.95%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'18.35"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -3.51%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'3.17"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -6.50%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  +0.30%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'0.119"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.55%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'4.33"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +10.42%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'21.92"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -2.48%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.08%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'3.32"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +2.15%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  -1.81%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'1.968.43"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -1.69%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'8.90"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -0.61%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = "'SEI"
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = "'https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = "'0.871"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +13.81%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'81.45"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -1.89%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = "'Algorand"
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'0.192"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.19%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'RocketPoolETH"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'2.735.19"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -2.48%  "
$ws.Range('E51').Style = 'Normal'
